$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value = "8e8f4ea7d0efe7065f4aa179c66a5993"
$ws.Range("B10").Value = "0df757e599079ee6e7287db471566dab"
$ws.Range("B58").Value = "e021118948136fc1197f1b99869af114"
$ws.Range("B71").Value = "7d5b1d4c9d76911c7f0629c2bbc3b559"
$ws.Range("B96").Value = "163a2c95fdc0133f3182e4a2f5981be1"
$ws.Range("B114").Value = "ee3e4de10c46cc607ae85f2e6657a31e"
$ws.Range("B120").Value = "8cef06adee08acc58c2564ba45a92776"
$ws.Range("B142").Value = "575e393b45ab9db58ab117dfedf0f70d"
$ws.Range("B190").Value = "a0e66fbb3a80f46243aa89c973e6aef5"
$ws.Range("B255").Value = "3c39cc40a5d3c996803a1bbb7835e95b"
$ws.Range("B352").Value = "444d7c36df66c5ffb38e38d0022965ff"
$ws.Range("B388").Value = "e021a1af0e663045acb12bbf52548523"
$ws.Range("B407").Value = "1eeacbd7d37f53f89db299ee668fff75"
$ws.Range("B419").Value = "afba4ee92bb44bede48ddf483ac24705"
$ws.Range("B472").Value = "846627bbd541c1508403cdd22739c10b"
$ws.Range("B492").Value = "ce84a2a5da4ea27b98021964a91beaa4"
$ws.Range("B500").Value = "59328d6fbee2ac587678815c09af1874"
$ws.Range("B561").Value = "c7bc39acd047929c20f71caa2141a1f2"
$ws.Range("B593").Value = "9b9367d22346d83cef61f20fb8cf1f46"
$ws.Range("B681").Value = "e1fa09aa78f53496969d261f9f5e7b69"
$ws.Range("B727").Value = "ea5085503eeecda17862f1fcddac8e01"
$ws.Range("B734").Value = "a885f747d9f8f8535cfd3087fce93e47"
$ws.Range("B776").Value = "e867a7ef5a2c4abfe453536e28f5ee67"
$ws.Range("B953").Value = "4f5e17e055f48fc2357151abfc4241f0"

$wb.Save()
